$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the _GoBack bookmark from its original location
# (right after "Team profile: " in the first Heading1 paragraph).
# It is later re-added at a different spot (see Change 3 below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: merge the runs around "follows" in the paragraph describing how
# a player may play treasure/action cards - i.e. drop the
# proofErr(gramStart)/"follows"/proofErr(gramEnd) wrapper while keeping the
# text identical. We only touch the portion of the paragraph AFTER the
# existing "deck,  A" proofErr pair (which must remain untouched).
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("deck,  A")
$afterAnchor = $anchor.End

$paraRng = $d.Content
$paraRng.Find.Execute("hand contains the cards drawn")
$paraEnd = $paraRng.Paragraphs(1).Range.End - 1   # exclude paragraph mark

$tailRange = $d.Range($afterAnchor, $paraEnd)
$tailText = $tailRange.Text
$tailRange.Delete()
$reinsert = $d.Range($afterAnchor, $afterAnchor)
$reinsert.InsertAfter($tailText)

# ---------------------------------------------------------------------------
# Change 3: add the _GoBack bookmark to the (empty) paragraph that sits
# right before "Name: Market" - i.e. right after the Laboratory card's
# "Comments:  " paragraph.
# ---------------------------------------------------------------------------
$marketRng = $d.Content
$marketRng.Find.Execute("Name: Market")
$marketPara = $marketRng.Paragraphs(1)
$targetPara = $marketPara.Previous()
$d.Bookmarks.Add("_GoBack", $targetPara.Range)

# ---------------------------------------------------------------------------
# Change 4a: fix grammar/wording in the "3. A.I.:" paragraph and drop its
# stale proofErr markers:
#   "decision making" -> "decision-making"
#   "well known"       -> "well-known"
#   "as  \u201c" (2 spaces) -> "as \u201c" (1 space)
# ---------------------------------------------------------------------------
$aiFind = $d.Content
$aiFind.Find.Execute("3. A.I.:")
$aiPara = $aiFind.Paragraphs(1)
$aiStart = $aiPara.Range.Start
$aiEnd = $aiPara.Range.End - 1   # exclude paragraph mark

$aiRange = $d.Range($aiStart, $aiEnd)
$aiOriginal = $aiRange.Text
$aiNew = $aiOriginal.Replace("decision making", "decision-making")
$aiNew = $aiNew.Replace("well known", "well-known")
$aiNew = $aiNew.Replace(("as  " + [char]8220), ("as " + [char]8220))

$aiRange.Delete()
$aiInsert = $d.Range($aiStart, $aiStart)
$aiInsert.InsertAfter($aiNew)

# ---------------------------------------------------------------------------
# Change 4b: in the "4. Gamelog:" paragraph, fix "game it" -> "game, it" and
# drop the proofErr(gramStart)/"game"/proofErr(gramEnd) wrapper, while
# leaving the other proofErr markers in that paragraph (around "Gamelog",
# "in order to", "gamelog") untouched.
# ---------------------------------------------------------------------------
$mpFind = $d.Content
$mpFind.Find.Execute("In a multiplayer ")
$gameStart = $mpFind.End
$gameEnd = $gameStart + 4   # length of "game"

# Extend one character on each side so the delete crosses both proofErr
# boundaries (they sit exactly at gameStart and gameEnd).
$gameRange = $d.Range($gameStart - 1, $gameEnd + 1)
$before = $d.Range($gameStart - 1, $gameStart).Text
$after = $d.Range($gameEnd, $gameEnd + 1).Text
$gameRange.Delete()
$gameInsert = $d.Range($gameStart - 1, $gameStart - 1)
$gameInsert.InsertAfter($before + "game," + $after)

Write-Output "edits applied"
